$d = $word.ActiveDocument

# 1. Fix typo: "começara" -> "começar"
$d.Content.Find.Execute("Alterar anos para começara no ano mais novo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Alterar anos para começar no ano mais novo", 2)

# 2. Locate the paragraph that now reads "Alterar anos para começar no ano mais novo"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Alterar anos para começar no ano mais novo*") {
        $target = $p
    }
}

# 3. Insert three new list paragraphs after it - they inherit the same list/paragraph
#    formatting (numId 2) from the paragraph being split. The last one gets a
#    trailing placeholder character "Z" so that the bookmark-move step below never
#    has to target the very last character position of a paragraph - a position the
#    bookmark engine can't reliably target directly.
$null = $target.Range.InsertParagraphAfter()
$p1 = $target.Next()
$null = $p1.Range.InsertAfter("Histórico")

$null = $p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$null = $p2.Range.InsertAfter("Carla - Erros")

$null = $p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$null = $p3.Range.InsertAfter("Unidades - DInheiroZ")

# 4. Move the _GoBack bookmark from the end of the (now fixed) "Alterar anos..."
#    paragraph to the end of the last new paragraph's real text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bmPos = $p3.Range.End - 2   # one position before the trailing "Z" placeholder
$bmRange = $d.Range($bmPos, $bmPos)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)

# 5. Strip the trailing "Z" placeholder back out, leaving the real target text with
#    the bookmark sitting right after it.
$zRange = $d.Range($p3.Range.End - 2, $p3.Range.End - 1)
$null = $zRange.Delete()
